$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Header date line
Replace-Text "2025-07-14 Monday" "2025-07-15 Tuesday"

# Row 1 (first block of 5)
Replace-Text "375÷7=" "675÷9="
Replace-Text "786÷8=" "761÷9="
Replace-Text "835÷3=" "300÷6="
Replace-Text "572÷3=" "966÷6="
Replace-Text "105÷8=" "244÷2="

# Row 2
Replace-Text "325÷6=" "776÷6="
Replace-Text "591÷9=" "504÷4="
Replace-Text "756÷4=" "337÷7="
Replace-Text "608÷5=" "345÷6="
Replace-Text "900÷7=" "553÷9="

# Row 3
Replace-Text "224÷9=" "707÷6="
Replace-Text "213÷7=" "356÷9="
Replace-Text "734÷9=" "736÷6="
Replace-Text "918÷7=" "889÷2="
Replace-Text "579÷9=" "597÷6="

# Row 4
Replace-Text "661÷3=" "786÷8="
Replace-Text "253÷6=" "701÷2="
Replace-Text "925÷6=" "637÷6="
Replace-Text "543÷6=" "650÷4="
Replace-Text "305÷3=" "103÷8="

# Row 5
Replace-Text "773÷3=" "623÷4="
Replace-Text "452÷8=" "698÷4="
Replace-Text "592÷3=" "525÷3="
Replace-Text "989÷3=" "116÷5="
Replace-Text "670÷7=" "327÷2="
